$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix formatting on fastq purpose column: "fullRNASEQ" -> "fullRNASeq"
for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    if ($cell.Text -eq "fullRNASEQ") {
        $cell.Value = "fullRNASeq"
    }
}
